$d = $word.ActiveDocument
$rng = $d.Content

# Step 1: Replace the header block (Provider Name / Addressed To / Reason for
# Letter / Illness lines) plus the "Dear ..." salutation and the four body
# paragraphs with the new letter body. Using MatchWildcards so the '*' can
# span the intervening <w:br/> line breaks, and '^l' in the replacement text
# inserts new manual line breaks at the paragraph-separator spots.
$found1 = $rng.Find.Execute(
    "Provider Name: Dr. John Smith*Thank you for your consideration. I look forward to hearing from you soon.",
    $false, $false, $true, $false, $false, $true, 1, $false,
    "Dear lhkbjb;knj,^l^lThis letter is being sent on behalf of jkn;jnj;l, a health provider practice. We are writing to inform you of a patient who has been diagnosed with pjh;lkj;.^l^lWe understand that this is a difficult time for the patient and their family, and we are here to provide the best care possible. Our team of experienced professionals is dedicated to providing the highest quality of care and support.^l^lWe are committed to providing the patient with the best possible treatment plan. We will work closely with the patient and their family to ensure that they receive the best care possible.^l^lWe understand that this is a difficult time for the patient and their family, and we are here to provide the best care possible. We are here to provide support and guidance throughout the entire process.^l^lThank you for your time and consideration. If you have any questions or concerns, please do not hesitate to contact us.",
    2)

# Step 2: Fix the signature block: "Sincerely,<br/><br/>Dr. John Smith"
# becomes "Sincerely,<br/>jkn;jnj;l" (one break removed, name swapped).
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "Sincerely,*Dr. John Smith",
    $false, $false, $true, $false, $false, $true, 1, $false,
    "Sincerely,^ljkn;jnj;l",
    2)

Write-Output "Block replace found: $found1; signature replace found: $found2"
